# "add auto remove function & refine the structure for the xlsx input"
#
# On sheet "爆款文案_2024" (the active sheet):
#   - Row 2 ("需要批量产生的句条数:" / 5) stays as-is but its label text is
#     refined (clarified) from "需要批量产生的条数:" to "需要批量产生的句条数:".
#   - A brand-new row 3 is introduced: "需要批量产生的文案条数:" / 10000,
#     formatted exactly like row 2 (same fill/border/number style).
#   - Column A is widened a bit to fit the longer labels.
#   - The remembered selection moves to A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- refine the existing label text on row 2 ---------------------------
$ws.Range("A2").Value = "需要批量产生的句条数:"

# --- insert the new "文案条数" row, cloning row 2's formatting ---------
$null = $ws.Range("A2:B2").Copy($ws.Range("A3:B3"))
$ws.Range("A3").Value = "需要批量产生的文案条数:"
$ws.Range("B3").Value = 10000

# --- widen column A to fit the new, longer labels -----------------------
$ws.Columns.Item(1).ColumnWidth = 24.5

# --- move the saved selection to A7 --------------------------------------
$null = $ws.Range("A7").Select()
